# Update the Five Year Plan workbook from FY2023-FY2027 to FY2024-FY2028
$wb = $excel.ActiveWorkbook

# ---------- Revenue by FY ----------
$ws = $wb.Worksheets.Item("Revenue by FY")
$ws.Range("A2").Value = "The Five Year Plan: FY 2024 - FY 2028"
$dateCell = $ws.Range("A4")
$origFmt = $dateCell.NumberFormat
$dateCell.NumberFormat = "@"
$dateCell.Value = "06/12/2023"
$dateCell.NumberFormat = $origFmt
$ws.Range("A6").Value = "FY 2024"
$ws.Range("A19").Value = "FY 2025"
$ws.Range("A32").Value = "FY 2026"
$ws.Range("A45").Value = "FY 2027"
$ws.Range("A58").Value = "FY 2028"
$ws.Range("A71").Value = "FY 2024 - FY 2028"

# ---------- Revenue by Tax ----------
$ws = $wb.Worksheets.Item("Revenue by Tax")
$ws.Range("A2").Value = "The Five Year Plan: FY 2024 - FY 2028"
$dateCell = $ws.Range("A4")
$origFmt = $dateCell.NumberFormat
$dateCell.NumberFormat = "@"
$dateCell.Value = "06/12/2023"
$dateCell.NumberFormat = $origFmt
$ws.Range("A7").Value = "FY 2024"
$ws.Range("A8").Value = "FY 2025"
$ws.Range("A9").Value = "FY 2026"
$ws.Range("A10").Value = "FY 2027"
$ws.Range("A11").Value = "FY 2028"
$ws.Range("A15").Value = "FY 2024"
$ws.Range("A16").Value = "FY 2025"
$ws.Range("A17").Value = "FY 2026"
$ws.Range("A18").Value = "FY 2027"
$ws.Range("A19").Value = "FY 2028"
$ws.Range("A23").Value = "FY 2024"
$ws.Range("A24").Value = "FY 2025"
$ws.Range("A25").Value = "FY 2026"
$ws.Range("A26").Value = "FY 2027"
$ws.Range("A27").Value = "FY 2028"
$ws.Range("A31").Value = "FY 2024"
$ws.Range("A32").Value = "FY 2025"
$ws.Range("A33").Value = "FY 2026"
$ws.Range("A34").Value = "FY 2027"
$ws.Range("A35").Value = "FY 2028"
$ws.Range("A39").Value = "FY 2024"
$ws.Range("A40").Value = "FY 2025"
$ws.Range("A41").Value = "FY 2026"
$ws.Range("A42").Value = "FY 2027"
$ws.Range("A43").Value = "FY 2028"
$ws.Range("A47").Value = "FY 2024"
$ws.Range("A48").Value = "FY 2025"
$ws.Range("A49").Value = "FY 2026"
$ws.Range("A50").Value = "FY 2027"
$ws.Range("A51").Value = "FY 2028"
$ws.Range("A55").Value = "FY 2024"
$ws.Range("A56").Value = "FY 2025"
$ws.Range("A57").Value = "FY 2026"
$ws.Range("A58").Value = "FY 2027"
$ws.Range("A59").Value = "FY 2028"
$ws.Range("A63").Value = "FY 2024"
$ws.Range("A64").Value = "FY 2025"
$ws.Range("A65").Value = "FY 2026"
$ws.Range("A66").Value = "FY 2027"
$ws.Range("A67").Value = "FY 2028"
$ws.Range("A71").Value = "FY 2024"
$ws.Range("A72").Value = "FY 2025"
$ws.Range("A73").Value = "FY 2026"
$ws.Range("A74").Value = "FY 2027"
$ws.Range("A75").Value = "FY 2028"
$ws.Range("A79").Value = "FY 2024"
$ws.Range("A80").Value = "FY 2025"
$ws.Range("A81").Value = "FY 2026"
$ws.Range("A82").Value = "FY 2027"
$ws.Range("A83").Value = "FY 2028"

# ---------- Growth by Tax ----------
$ws = $wb.Worksheets.Item("Growth by Tax")
$ws.Range("B2").Value = "The Five Year Plan: FY 2024 - FY 2028"
$ws.Range("B9").Value = "FY 2024"
$ws.Range("B10").Value = "FY 2025"
$ws.Range("B11").Value = "FY 2026"
$ws.Range("B12").Value = "FY 2027"
$ws.Range("B13").Value = "FY 2028"
$ws.Range("B16").Value = "FY 2024"
$ws.Range("B17").Value = "FY 2025"
$ws.Range("B18").Value = "FY 2026"
$ws.Range("B19").Value = "FY 2027"
$ws.Range("B20").Value = "FY 2028"
$ws.Range("B23").Value = "FY 2024"
$ws.Range("B24").Value = "FY 2025"
$ws.Range("B25").Value = "FY 2026"
$ws.Range("B26").Value = "FY 2027"
$ws.Range("B27").Value = "FY 2028"
$ws.Range("B30").Value = "FY 2024"
$ws.Range("B31").Value = "FY 2025"
$ws.Range("B32").Value = "FY 2026"
$ws.Range("B33").Value = "FY 2027"
$ws.Range("B34").Value = "FY 2028"
$ws.Range("B37").Value = "FY 2024"
$ws.Range("B38").Value = "FY 2025"
$ws.Range("B39").Value = "FY 2026"
$ws.Range("B40").Value = "FY 2027"
$ws.Range("B41").Value = "FY 2028"
$ws.Range("B44").Value = "FY 2024"
$ws.Range("B45").Value = "FY 2025"
$ws.Range("B46").Value = "FY 2026"
$ws.Range("B47").Value = "FY 2027"
$ws.Range("B48").Value = "FY 2028"
$ws.Range("B51").Value = "FY 2024"
$ws.Range("B52").Value = "FY 2025"
$ws.Range("B53").Value = "FY 2026"
$ws.Range("B54").Value = "FY 2027"
$ws.Range("B55").Value = "FY 2028"
$ws.Range("B58").Value = "FY 2024"
$ws.Range("B59").Value = "FY 2025"
$ws.Range("B60").Value = "FY 2026"
$ws.Range("B61").Value = "FY 2027"
$ws.Range("B62").Value = "FY 2028"

# ---------- Growth by Year ----------
$ws = $wb.Worksheets.Item("Growth by Year")
$ws.Range("B2").Value = "The Five Year Plan: FY 2024 - FY 2028"
$ws.Range("B7").Value = "FY 2024"
$ws.Range("B17").Value = "FY 2025"
$ws.Range("B27").Value = "FY 2026"
$ws.Range("B37").Value = "FY 2027"
$ws.Range("B47").Value = "FY 2028"

# ---------- Revenue Data ----------
$ws = $wb.Worksheets.Item("Revenue Data")
$ws.Cells.Item(2, 1).Value = 2023
$ws.Cells.Item(2, 2).Value = 1703381000
$ws.Cells.Item(2, 3).Value = 1703381000
$ws.Cells.Item(3, 1).Value = 2024
$ws.Cells.Item(3, 2).Value = 1753427000
$ws.Cells.Item(3, 3).Value = 1772269627.847429
$ws.Cells.Item(4, 1).Value = 2025
$ws.Cells.Item(4, 2).Value = 1836714000
$ws.Cells.Item(4, 3).Value = 1838334302.577886
$ws.Cells.Item(5, 1).Value = 2026
$ws.Cells.Item(5, 2).Value = 1920285000
$ws.Cells.Item(5, 3).Value = 1918285780.126404
$ws.Cells.Item(6, 1).Value = 2027
$ws.Cells.Item(6, 2).Value = 1998248000
$ws.Cells.Item(6, 3).Value = 1997262421.65557
$ws.Cells.Item(7, 1).Value = 2028
$ws.Cells.Item(7, 2).Value = 2075780000
$ws.Cells.Item(7, 3).Value = 2079080418.530375
$ws.Cells.Item(8, 1).Value = 2023
$ws.Cells.Item(8, 2).Value = 401602641.4853873
$ws.Cells.Item(8, 3).Value = 401602641.4853873
$ws.Cells.Item(9, 1).Value = 2024
$ws.Cells.Item(9, 2).Value = 417465946
$ws.Cells.Item(9, 3).Value = 426102781.8896299
$ws.Cells.Item(10, 1).Value = 2025
$ws.Cells.Item(10, 2).Value = 434957769
$ws.Cells.Item(10, 3).Value = 439712662.813524
$ws.Cells.Item(11, 1).Value = 2026
$ws.Cells.Item(11, 2).Value = 451442668
$ws.Cells.Item(11, 3).Value = 451337933.3409257
$ws.Cells.Item(12, 1).Value = 2027
$ws.Cells.Item(12, 2).Value = 467468883
$ws.Cells.Item(12, 3).Value = 463871970.5891466
$ws.Cells.Item(13, 1).Value = 2028
$ws.Cells.Item(13, 2).Value = 483646001
$ws.Cells.Item(13, 3).Value = 475938749.0414606
$ws.Cells.Item(14, 1).Value = 2023
$ws.Cells.Item(14, 2).Value = 729455000
$ws.Cells.Item(14, 3).Value = 729455000
$ws.Cells.Item(15, 1).Value = 2024
$ws.Cells.Item(15, 2).Value = 709051000
$ws.Cells.Item(15, 3).Value = 719482885.7161967
$ws.Cells.Item(16, 1).Value = 2025
$ws.Cells.Item(16, 2).Value = 725688000
$ws.Cells.Item(16, 3).Value = 732518837.8059332
$ws.Cells.Item(17, 1).Value = 2026
$ws.Cells.Item(17, 2).Value = 708707000
$ws.Cells.Item(17, 3).Value = 750774788.5591393
$ws.Cells.Item(18, 1).Value = 2027
$ws.Cells.Item(18, 2).Value = 757182000
$ws.Cells.Item(18, 3).Value = 767210127.5964483
$ws.Cells.Item(19, 1).Value = 2028
$ws.Cells.Item(19, 2).Value = 772629000
$ws.Cells.Item(19, 3).Value = 783119205.9714768
$ws.Cells.Item(20, 1).Value = 2023
$ws.Cells.Item(20, 2).Value = 394167835.5139231
$ws.Cells.Item(20, 3).Value = 394167835.5139231
$ws.Cells.Item(21, 1).Value = 2024
$ws.Cells.Item(21, 2).Value = 432300000
$ws.Cells.Item(21, 3).Value = 333285196.13617
$ws.Cells.Item(22, 1).Value = 2025
$ws.Cells.Item(22, 2).Value = 440946000
$ws.Cells.Item(22, 3).Value = 337661362.6652149
$ws.Cells.Item(23, 1).Value = 2026
$ws.Cells.Item(23, 2).Value = 449765000
$ws.Cells.Item(23, 3).Value = 387702614.4153731
$ws.Cells.Item(24, 1).Value = 2027
$ws.Cells.Item(24, 2).Value = 458760000
$ws.Cells.Item(24, 3).Value = 446573572.9402359
$ws.Cells.Item(25, 1).Value = 2028
$ws.Cells.Item(25, 2).Value = 467981000
$ws.Cells.Item(25, 3).Value = 482886273.8934793
$ws.Cells.Item(26, 1).Value = 2023
$ws.Cells.Item(26, 2).Value = 97162999.99999993
$ws.Cells.Item(26, 3).Value = 97162999.99999993
$ws.Cells.Item(27, 1).Value = 2024
$ws.Cells.Item(27, 2).Value = 100087000
$ws.Cells.Item(27, 3).Value = 101351353.625347
$ws.Cells.Item(28, 1).Value = 2025
$ws.Cells.Item(28, 2).Value = 103841000
$ws.Cells.Item(28, 3).Value = 103859752.856211
$ws.Cells.Item(29, 1).Value = 2026
$ws.Cells.Item(29, 2).Value = 108192000
$ws.Cells.Item(29, 3).Value = 107246318.6741039
$ws.Cells.Item(30, 1).Value = 2027
$ws.Cells.Item(30, 2).Value = 112379000
$ws.Cells.Item(30, 3).Value = 111142611.2376092
$ws.Cells.Item(31, 1).Value = 2028
$ws.Cells.Item(31, 2).Value = 116750000
$ws.Cells.Item(31, 3).Value = 115033780.2534414
$ws.Cells.Item(32, 1).Value = 2023
$ws.Cells.Item(32, 2).Value = 33128000
$ws.Cells.Item(32, 3).Value = 33128000
$ws.Cells.Item(33, 1).Value = 2024
$ws.Cells.Item(33, 2).Value = 38379000
$ws.Cells.Item(33, 3).Value = 34287580.11504535
$ws.Cells.Item(34, 1).Value = 2025
$ws.Cells.Item(34, 2).Value = 39404000
$ws.Cells.Item(34, 3).Value = 36214252.59146186
$ws.Cells.Item(35, 1).Value = 2026
$ws.Cells.Item(35, 2).Value = 40283000
$ws.Cells.Item(35, 3).Value = 38011180.79252171
$ws.Cells.Item(36, 1).Value = 2027
$ws.Cells.Item(36, 2).Value = 41109000
$ws.Cells.Item(36, 3).Value = 39870409.2333933
$ws.Cells.Item(37, 1).Value = 2028
$ws.Cells.Item(37, 2).Value = 41948000
$ws.Cells.Item(37, 3).Value = 41818825.39507025
$ws.Cells.Item(38, 1).Value = 2023
$ws.Cells.Item(38, 2).Value = 29896000
$ws.Cells.Item(38, 3).Value = 29896000
$ws.Cells.Item(39, 1).Value = 2024
$ws.Cells.Item(39, 2).Value = 30944000
$ws.Cells.Item(39, 3).Value = 27757257.39004321
$ws.Cells.Item(40, 1).Value = 2025
$ws.Cells.Item(40, 2).Value = 32779000
$ws.Cells.Item(40, 3).Value = 29482658.77861664
$ws.Cells.Item(41, 1).Value = 2026
$ws.Cells.Item(41, 2).Value = 33687000
$ws.Cells.Item(41, 3).Value = 29952111.02690598
$ws.Cells.Item(42, 1).Value = 2027
$ws.Cells.Item(42, 2).Value = 35539000
$ws.Cells.Item(42, 3).Value = 30742392.88784411
$ws.Cells.Item(43, 1).Value = 2028
$ws.Cells.Item(43, 2).Value = 37697000
$ws.Cells.Item(43, 3).Value = 31561717.13324906

# ---------- Tax Base Data ----------
$ws = $wb.Worksheets.Item("Tax Base Data")
$ws.Cells.Item(2, 1).Value = 44927
$ws.Cells.Item(2, 2).Value = 63318006096.20104
$ws.Cells.Item(2, 3).Value = 20080132074.26936
$ws.Cells.Item(2, 4).Value = 20997590908.28648
$ws.Cells.Item(2, 5).Value = 11547556183.74558
$ws.Cells.Item(2, 6).Value = 9450034724.540899
$ws.Cells.Item(2, 7).Value = 12024644158.44793
$ws.Cells.Item(2, 8).Value = 431835555.5555552
$ws.Cells.Item(2, 9).Value = 662560000
$ws.Cells.Item(2, 10).Value = 1120623734.912662
$ws.Cells.Item(3, 1).Value = 45292
$ws.Cells.Item(3, 2).Value = 66417985437.01514
$ws.Cells.Item(3, 3).Value = 21305139094.48149
$ws.Cells.Item(3, 4).Value = 21570020833.13495
$ws.Cells.Item(3, 5).Value = 12186847765.69809
$ws.Cells.Item(3, 6).Value = 9383173067.436855
$ws.Cells.Item(3, 7).Value = 10167333622.21385
$ws.Cells.Item(3, 8).Value = 450450460.5570976
$ws.Cells.Item(3, 9).Value = 685751602.300907
$ws.Cells.Item(3, 10).Value = 1046605600.724371
$ws.Cells.Item(4, 1).Value = 45658
$ws.Cells.Item(4, 2).Value = 68893840428.37903
$ws.Cells.Item(4, 3).Value = 21985633140.6762
$ws.Cells.Item(4, 4).Value = 22008697271.55444
$ws.Cells.Item(4, 5).Value = 12470854204.43241
$ws.Cells.Item(4, 6).Value = 9537843067.122036
$ws.Cells.Item(4, 7).Value = 10300834736.58374
$ws.Cells.Item(4, 8).Value = 461598901.5831602
$ws.Cells.Item(4, 9).Value = 724285051.8292372
$ws.Cells.Item(4, 10).Value = 1111662992.072636
$ws.Cells.Item(5, 1).Value = 46023
$ws.Cells.Item(5, 2).Value = 71890120445.84735
$ws.Cells.Item(5, 3).Value = 22566896667.04628
$ws.Cells.Item(5, 4).Value = 22691130588.40522
$ws.Cells.Item(5, 5).Value = 12958507921.74145
$ws.Cells.Item(5, 6).Value = 9732622666.663769
$ws.Cells.Item(5, 7).Value = 11827413496.50314
$ws.Cells.Item(5, 8).Value = 476650305.2182394
$ws.Cells.Item(5, 9).Value = 760223615.8504341
$ws.Cells.Item(5, 10).Value = 1129363997.090103
$ws.Cells.Item(6, 1).Value = 46388
$ws.Cells.Item(6, 2).Value = 74849867283.75912
$ws.Cells.Item(6, 3).Value = 23193598529.45733
$ws.Cells.Item(6, 4).Value = 23231158787.86824
$ws.Cells.Item(6, 5).Value = 13299352881.90872
$ws.Cells.Item(6, 6).Value = 9931805905.959518
$ws.Cells.Item(6, 7).Value = 13623354879.20182
$ws.Cells.Item(6, 8).Value = 493967161.0560406
$ws.Cells.Item(6, 9).Value = 797408184.667866
$ws.Cells.Item(6, 10).Value = 1159162093.140666
$ws.Cells.Item(7, 1).Value = 46753
$ws.Cells.Item(7, 2).Value = 77916097410.1067
$ws.Cells.Item(7, 3).Value = 23796937452.07303
$ws.Cells.Item(7, 4).Value = 23755142007.57355
$ws.Cells.Item(7, 5).Value = 13630930307.36266
$ws.Cells.Item(7, 6).Value = 10124211700.21089
$ws.Cells.Item(7, 7).Value = 14731124889.97801
$ws.Cells.Item(7, 8).Value = 511261245.5708506
$ws.Cells.Item(7, 9).Value = 836376507.9014051
$ws.Cells.Item(7, 10).Value = 1190055251.351524
